$wb = $excel.ActiveWorkbook

# The new question lives on the 15th worksheet (index 15, 1-based) -> sheet name "14_"
$ws = $wb.Worksheets.Item(15)

# --- Cell content -----------------------------------------------------
# NB: new shared-string entries are appended in first-write order, so the
# question/answer text is written before the explanation text in C2 to
# reproduce the same shared-string ordering as the authored workbook.
$ws.Range("A1").Value2 = "Consider a simple circuit with a battery and a resistor.   If the resistor is made smaller, what happens to the current?"

$ws.Range("A2").Value2 = "It increases"
$ws.Range("A3").Value2 = "It stays the same"
$ws.Range("A4").Value2 = "It decreases"
$ws.Range("A5").Value2 = "It depends on the voltage of the battery"

$ws.Range("B2").Value2 = "Y"
$ws.Range("B3").Value2 = "N"
$ws.Range("B4").Value2 = "N"
$ws.Range("B5").Value2 = "N"

$ws.Range("C2").Value2 = "Yep!  If the driving force for a flow (voltage) stays the same and the resistance to the flow decreases, the flow (i.e. the current) will increase."

# --- Formatting ---------------------------------------------------------
# Apply the same wrap-text style used by every other cell on this sheet
# across the whole used block (A1:E13), matching the neighbouring sheets.
$ws.Range("A1:E13").WrapText = $true

# Row heights (points) that fit the wrapped question/answer text.
$ws.Rows.Item(1).RowHeight = 75
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 30

# Column widths for the question/explanation columns (closest values the
# host's character-width rounding can reproduce for the authored
# 23.5703125 / 29.42578125 pixel-derived widths).
$ws.Columns.Item(1).ColumnWidth = 22.666666666666664
$ws.Columns.Item(3).ColumnWidth = 28.666666666666664

# --- Sheet/view selection -------------------------------------------------
# Make this newly-populated sheet the active tab, and select G7 on it
# (this also clears tabSelected on whichever sheet previously had it).
$ws.Activate()
$ws.Range("G7").Select()
